# Apply the commit:
#  - "sua chiet khau cua sale phu"  -> adjust AC (Day du) / AF (Tong cong)
#    numbers for a handful of rows
#  - "update chien luoc chay tinh luong theo gio" -> bump last_edited_time
#    (column D) for every data row to the new run timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: last_edited_time -> 2024-07-21T16:44:00.000Z for rows 2..20
$newTimestamp = "2024-07-21T16:44:00.000Z"
for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 4).Value = $newTimestamp
}

# --- Columns AC (Day du) / AF (Tong cong): updated numeric values
$ws.Cells.Item(8, 29).Value  = 20     # AC8
$ws.Cells.Item(8, 32).Value  = 21     # AF8

$ws.Cells.Item(11, 29).Value = 20     # AC11
$ws.Cells.Item(11, 32).Value = 22     # AF11

$ws.Cells.Item(14, 29).Value = 20     # AC14
$ws.Cells.Item(14, 32).Value = 20.5   # AF14

$ws.Cells.Item(17, 29).Value = 20     # AC17
$ws.Cells.Item(17, 32).Value = 20.5   # AF17

$ws.Cells.Item(19, 29).Value = 20     # AC19
$ws.Cells.Item(19, 32).Value = 21     # AF19

$ws.Cells.Item(20, 29).Value = 21     # AC20
$ws.Cells.Item(20, 32).Value = 21     # AF20
